$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.085.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.95%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.888.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.65%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9989"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.91%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9984"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5188"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.04%  "

$ws.Range("E8").Value = "  +2.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07211"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.67%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9048"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.81%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.08"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07664"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.81%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "95.23"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.868.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.270"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.68%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9989"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.16%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008507"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.33%  "

$ws.Range("E18").Value = "  +1.89%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9987"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.127.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.94%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.049"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.75%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.109.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.461"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.29%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.35%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.787"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.87%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.131"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.23%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.933"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.59%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.795"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09208"
$ws.Range("D32").Style = "Normal"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05044"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.98%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7628"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.89%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.192"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.79%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.018"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.55%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.281"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.74%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.549"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5611"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.66%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01994"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.53%  "

$ws.Range("E41").Value = "  +0.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.595"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "118.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.857"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1508"
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4793"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.51%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.22%  "

$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9987"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.08%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.576"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.13"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.38%  "

$ws.Range("E51").Value = "  +0.96%  "
